$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "AA6" = 17
    "AB6" = 51
    "AC6" = 51
    "AD6" = 67
    "AE6" = 5.5
    "AF6" = 6
    "AG6" = 23
    "AH6" = 101
    "AJ6" = 4.75
    "AK6" = 7.5
    "AP6" = 2.1
    "AQ6" = 1.78
    "AR6" = 5
    "AS6" = 1.17
    "G6" = 4.75
    "H6" = 3
    "I6" = 1.95
    "J6" = 5.5
    "K6" = 1.83
    "L6" = 2.75
    "M6" = 1.13
    "N6" = 6
    "O6" = 1.62
    "P6" = 2.2
    "Q6" = 2.88
    "R6" = 1.4
    "S6" = 6.5
    "T6" = 1.11
    "U6" = 1.67
    "V6" = 2.1
    "W6" = 2.5
    "X6" = 1.5
    "Y6" = 8.5
    "Z6" = 21
    "AA7" = 10
    "AB7" = 17
    "AE7" = 6
    "AF7" = 6
    "AG7" = 21
    "AH7" = 81
    "AJ7" = 9
    "AK7" = 21
    "AL7" = 17
    "AN7" = 41
    "AO7" = 51
    "AR7" = 4.9
    "AS7" = 1.18
    "G7" = 1.95
    "H7" = 3
    "I7" = 4.75
    "J7" = 2.75
    "L7" = 5.5
    "M7" = 1.13
    "N7" = 6
    "W7" = 2.38
    "X7" = 1.53
    "Y7" = 5
    "Z7" = 7.5
    "AA8" = 10
    "AB8" = 26
    "AC8" = 21
    "AD8" = 29
    "AE8" = 10
    "AG8" = 15
    "AI8" = 251
    "AJ8" = 8.5
    "AK8" = 13
    "AL8" = 10
    "AM8" = 26
    "AO8" = 29
    "AR8" = 2.95
    "AS8" = 1.41
    "G8" = 2.6
    "H8" = 3.5
    "I8" = 2.63
    "J8" = 3.25
    "K8" = 2.1
    "L8" = 3.4
    "M8" = 1.06
    "N8" = 10
    "O8" = 1.33
    "P8" = 3.25
    "Q8" = 2.05
    "R8" = 1.75
    "S8" = 3.5
    "T8" = 1.29
    "W8" = 1.8
    "X8" = 1.95
    "Z8" = 12
    "AA18" = 9
    "AB18" = 19
    "AC18" = 19
    "AD18" = 29
    "AE18" = 8.5
    "AG18" = 15
    "AH18" = 51
    "AI18" = 301
    "AJ18" = 9.5
    "AK18" = 17
    "AM18" = 41
    "AP18" = 1.8
    "AQ18" = 2.05
    "AR18" = 3
    "AS18" = 1.37
    "G18" = 2.1
    "H18" = 3.25
    "I18" = 3.5
    "J18" = 2.88
    "K18" = 2.05
    "M18" = 1.06
    "N18" = 10
    "O18" = 1.3
    "P18" = 3.4
    "Q18" = 2.05
    "R18" = 1.75
    "S18" = 3.5
    "T18" = 1.29
    "U18" = 1.44
    "V18" = 2.63
    "W18" = 1.83
    "X18" = 1.83
    "Y18" = 7
    "Z18" = 9.5
    "AC19" = 29
    "AD19" = 51
    "AE19" = 5
    "M19" = 1.14
    "O19" = 1.67
    "T19" = 1.1
    "U19" = 1.73
    "V19" = 2
    "Q25" = 2.1
    "R25" = 1.7
    "L51" = 3.6
    "M51" = 1.05
    "N51" = 11
    "O51" = 1.25
    "P51" = 3.75
    "Q51" = 1.88
    "R51" = 1.98
    "S51" = 3
    "T51" = 1.36
    "W51" = 1.67
    "X51" = 2.1
    "AK52" = 41
    "AL52" = 21
    "G52" = 1.4
    "I52" = 7
    "J52" = 1.91
    "L52" = 7
    "M52" = 1.03
    "N52" = 15
    "AA53" = 13
    "AE53" = 6
    "AJ53" = 5.5
    "AK53" = 10
    "AM53" = 23
    "AN53" = 26
    "G53" = 3.1
    "H53" = 2.8
    "I53" = 2.3
    "J53" = 4.33
    "L53" = 3.4
    "O53" = 1.62
    "P53" = 2.2
    "S53" = 6.5
    "T53" = 1.11
    "U53" = 1.62
    "V53" = 2.2
    "Y53" = 7
    "Z53" = 15
    "AG54" = 26
    "M54" = 1.05
    "N54" = 11
    "W54" = 2.25
    "X54" = 1.57
    "Y54" = 6.5
    "AP66" = 1.25
    "AQ66" = 4
    "AB93" = 12.5
    "AH93" = 120
    "AJ93" = 10.75
    "AK93" = 25
    "AL93" = 16
    "AM93" = 80
    "G93" = 1.7
    "I93" = 4.6
    "J93" = 2.3
    "L93" = 5
    "P93" = 2.67
    "Q93" = 2.05
    "R93" = 1.6
    "T93" = 1.22
    "V93" = 2.37
    "W93" = 2.02
    "X93" = 1.62
    "Z93" = 7
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Updated $($updates.Count) cells"